$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibition)
$wsExhibition = $wb.Worksheets.Item("展览")

$wsExhibition.Range("F3").Value = 7404
$wsExhibition.Range("F4").Value = 3532
$wsExhibition.Range("F6").Value = 3862
$wsExhibition.Range("F7").Value = 69
$wsExhibition.Range("F8").Value = 88
$wsExhibition.Range("F10").Value = 104
$wsExhibition.Range("F11").Value = 160
$wsExhibition.Range("F12").Value = 513
$wsExhibition.Range("F18").Value = 353
$wsExhibition.Range("F19").Value = 4159
$wsExhibition.Range("F24").Value = 1875
$wsExhibition.Range("G24").Value = 39.9
$wsExhibition.Range("F25").Value = 117
$wsExhibition.Range("F27").Value = 3065
$wsExhibition.Range("F28").Value = 2281
$wsExhibition.Range("F29").Value = 66
$wsExhibition.Range("F30").Value = 84
$wsExhibition.Range("F32").Value = 42
$wsExhibition.Range("F36").Value = 4357
$wsExhibition.Range("F37").Value = 490
$wsExhibition.Range("F41").Value = 824
$wsExhibition.Range("F42").Value = 224
$wsExhibition.Range("F44").Value = 1651
$wsExhibition.Range("F47").Value = 615

# Sheet 4: 全部类型 (All Types)
$wsAll = $wb.Worksheets.Item("全部类型")

$wsAll.Range("F5").Value = 7405
$wsAll.Range("F6").Value = 3532
$wsAll.Range("F7").Value = 3862
$wsAll.Range("F8").Value = 69
$wsAll.Range("F9").Value = 88
$wsAll.Range("F11").Value = 104
$wsAll.Range("F13").Value = 160
$wsAll.Range("F14").Value = 513
$wsAll.Range("F20").Value = 353
$wsAll.Range("F21").Value = 4159
$wsAll.Range("F27").Value = 1875
$wsAll.Range("G27").Value = 39.9
$wsAll.Range("F28").Value = 117
$wsAll.Range("F30").Value = 3065
$wsAll.Range("F31").Value = 2281
$wsAll.Range("F32").Value = 66
$wsAll.Range("F33").Value = 84
$wsAll.Range("F37").Value = 4357
$wsAll.Range("F39").Value = 490
$wsAll.Range("F42").Value = 824
$wsAll.Range("F43").Value = 224
$wsAll.Range("F45").Value = 1651
$wsAll.Range("F48").Value = 615

Write-Output "Update complete"
